# The source workbook (confirmations/102_1.xlsx) got two "Withdrawn" tallies
# filled in that were previously left blank:
#   B17 -> "     Civilian (FS, PHS, CG, NOAA), Withdrawn "  = 1
#   B22 -> "     Air Force, Withdrawn "                     = 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = 1
$ws.Range("B22").Value = 1
